# Updated cryptos list on Thu Jun 22 16:53:44 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.910.76'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.876.58'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.88'
$ws.Range('E5').Value = '  -2.93%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4924'
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.09'
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2900'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06586'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').Value = '1.874.69'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.88'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07178'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6652'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '85.18'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.828'
$ws.Range('E16').Value = '  +0.44%  '
$ws.Range('D17').Value = '29.913.70'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007819'
$ws.Range('E18').Value = '  +3.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9977'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.75'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').Value = '2.119.47'
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.746'
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.561'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.094'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '147.82'
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '134.50'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.69'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.914'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.380'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.168'
$ws.Range('E31').Value = '  -1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08595'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.914'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04979'
$ws.Range('E34').Value = '  -1.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.104'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7028'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.656'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.196'
$ws.Range('E38').Value = '  -5.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.681'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9301'
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01636'
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.032'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9946'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.55'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4155'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.579'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1254'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05707'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '32.55'
$ws.Range('E49').Value = '  +0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.198'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.331'
$ws.Range('E51').Value = '  -0.89%  '
